$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '69.686.15'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = '3.486.92'
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '192.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.29%  '
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.211'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.658'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.32'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000305'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.57'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.65%  '
$ws.Range('D14').Value = '4.053.88'
$ws.Range('E14').Value = '  -1.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '603.60'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.34%  '
$ws.Range('D16').Value = '69.842.61'
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '12.63'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.80'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('D19').Value = '3.483.86'
$ws.Range('E19').Value = '  -2.00%  '
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '105.35'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +11.24%  '
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.55%  '
$ws.Range('E27').Value = '  -0.57%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.84'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.09'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +13.18%  '
$ws.Range('E32').Value = '  +3.54%  '
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '64.21'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.93%  '
$ws.Range('D35').Value = '3.702.36'
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('E37').Value = '  -5.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '517.58'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('D39').Value = '0.0₃0790'
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.58'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.86%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.390'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.51'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.67%  '
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0461'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.12%  '
$ws.Range('E46').Value = '  +1.61%  '
$ws.Range('E47').Value = '  -3.92%  '
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.72'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.69'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.30'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.87%  '
